$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.981300392979479
$ws.Range("E2").Value = 2.228397109637226

$ws.Range("C3").Value = -0.4626567965509865
$ws.Range("E3").Value = -0.2018858887078534

$ws.Range("C4").Value = -0.4861108058567654
$ws.Range("E4").Value = 0.2108047537406676

$ws.Range("C5").Value = 2.490556233265284
$ws.Range("E5").Value = 0.9684279156219944

$ws.Range("C6").Value = 1.175152734113438
$ws.Range("E6").Value = 1.69818237209749

$ws.Range("C7").Value = -0.2674335569108899
$ws.Range("E7").Value = 0.679606031449409

$ws.Range("C8").Value = 2.038609866767938
$ws.Range("E8").Value = 1.325176859452326

$ws.Range("C9").Value = 1.566972126803345
$ws.Range("E9").Value = 1.522808462763714

$ws.Range("C10").Value = 2.246337373619012
$ws.Range("E10").Value = 1.693557061600992

$ws.Range("C11").Value = 1.777150434343522
$ws.Range("E11").Value = 1.905564797014625

$ws.Range("C12").Value = 1.741137453897301
$ws.Range("E12").Value = 2.082477074609068

$ws.Range("C13").Value = 1.562095320687407
$ws.Range("E13").Value = 1.845103901518885

$ws.Range("C14").Value = -0.8985012482809474
$ws.Range("E14").Value = 0.02570757229449772

$ws.Range("C15").Value = -0.7241284594088016
$ws.Range("E15").Value = -0.1006764062950749

$ws.Range("C16").Value = 3.72986100763808
$ws.Range("E16").Value = 1.357201584213352

$ws.Range("C17").Value = -0.3358560228777674
$ws.Range("E17").Value = 1.155533515694374

$ws.Range("C18").Value = -1.054811008161194
$ws.Range("E18").Value = -0.009717596728553435

$ws.Range("C19").Value = 1.593155398714674
$ws.Range("E19").Value = 0.1924237296221154
